$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F3 86 -> 87, F4 980 -> 990
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 87
$ws1.Range("F4").Value = 990

# Sheet "全部类型" (sheet4): update F3 86 -> 87, F4 980 -> 990
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 87
$ws4.Range("F4").Value = 990
